# Insert a new "Author" paragraph for the affiliation, right after the
# existing "Edison Achalma" author-name paragraph (inside the title block).

$d = $word.ActiveDocument

# Locate the target paragraph: the one whose text is exactly
# "Edison Achalma" and which uses the "Author" paragraph style.
$targetIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Edison Achalma" -and $p.Style.NameLocal -eq "Author") {
        $targetIndex = $i
    }
    $i = $i + 1
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'Edison Achalma' Author paragraph"
}

$authorParagraph = $d.Paragraphs.Item($targetIndex + 1)

# Insert a new empty paragraph right after it.
$endRange = $d.Range($authorParagraph.Range.End, $authorParagraph.Range.End)
$endRange.InsertParagraphAfter()

# The freshly inserted paragraph is now the next one in the collection;
# give it the "Author" style and fill in the affiliation text.
$newParagraph = $d.Paragraphs.Item($targetIndex + 2)
$newParagraph.Style = "Author"
$newParagraph.Range.Text = "Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga"
